$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cells (rows with changed C/D/E/F values) ---
# Row 4
$ws.Range("C4").Value = -1
$ws.Range("D4").Value = 45833.77620284099
$ws.Range("E4").Value = -1
$ws.Range("F4").Value = 45833.69550925926
# Row 5
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 45833.77620283774
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = 45833.59693287037
# Row 12
$ws.Range("C12").Value = 60
$ws.Range("D12").Value = 45833.77620284134
$ws.Range("E12").Value = 60
$ws.Range("F12").Value = 45833.73585648148
# Row 24
$ws.Range("C24").Value = 47
$ws.Range("D24").Value = 45833.77618177847
$ws.Range("E24").Value = 47
$ws.Range("F24").Value = 45833.54157407407
# Row 33
$ws.Range("C33").Value = 2831
$ws.Range("D33").Value = 45833.77620283895
$ws.Range("E33").Value = 2831
$ws.Range("F33").Value = 45833.68387731481
# Row 56
$ws.Range("C56").Value = 297
$ws.Range("D56").Value = 45833.77620284142
$ws.Range("E56").Value = 297
$ws.Range("F56").Value = 45833.73585648148
# Row 58
$ws.Range("C58").Value = 152
$ws.Range("D58").Value = 45833.77620284149
$ws.Range("E58").Value = 152
$ws.Range("F58").Value = 45833.73585648148
# Row 59
$ws.Range("C59").Value = 123
$ws.Range("D59").Value = 45833.77618177872
$ws.Range("E59").Value = 123
$ws.Range("F59").Value = 45833.54157407407
# Row 61
$ws.Range("C61").Value = 34
$ws.Range("D61").Value = 45834.28681886554
$ws.Range("E61").Value = 34
$ws.Range("F61").Value = 45833.87444444445
# Row 63
$ws.Range("C63").Value = 140
$ws.Range("D63").Value = 45833.77620283903
$ws.Range("E63").Value = 140
$ws.Range("F63").Value = 45833.68387731481
# Row 69
$ws.Range("C69").Value = 33
$ws.Range("D69").Value = 45833.77618177396
$ws.Range("E69").Value = 33
$ws.Range("F69").Value = 45833.47572916667
# Row 70
$ws.Range("C70").Value = 81
$ws.Range("D70").Value = 45833.77618177376
$ws.Range("E70").Value = 81
$ws.Range("F70").Value = 45833.47356481481
# Row 81
$ws.Range("C81").Value = 199
$ws.Range("D81").Value = 45833.77620283993
$ws.Range("E81").Value = 199
$ws.Range("F81").Value = 45833.68471064815
# Row 83
$ws.Range("C83").Value = 252
$ws.Range("D83").Value = 45833.77618176836
$ws.Range("E83").Value = 252
$ws.Range("F83").Value = 45833.39194444445
# Row 87
$ws.Range("C87").Value = 98
$ws.Range("D87").Value = 45833.77618176855
$ws.Range("E87").Value = 98
$ws.Range("F87").Value = 45833.39194444445
# Row 91
$ws.Range("C91").Value = 88
$ws.Range("D91").Value = 45833.77618176876
$ws.Range("E91").Value = 88
$ws.Range("F91").Value = 45833.39194444445
# Row 94
$ws.Range("C94").Value = 132
$ws.Range("D94").Value = 45833.77620284
$ws.Range("E94").Value = 132
$ws.Range("F94").Value = 45833.68471064815
# Row 101
$ws.Range("C101").Value = 2796
$ws.Range("D101").Value = 45833.77620284156
$ws.Range("E101").Value = 2796
$ws.Range("F101").Value = 45833.73585648148
# Row 117
$ws.Range("C117").Value = 1011
$ws.Range("D117").Value = 45833.7761817699
$ws.Range("E117").Value = 1011
$ws.Range("F117").Value = 45833.42523148148
# Row 121
$ws.Range("C121").Value = 638
$ws.Range("D121").Value = 45833.7761817701
$ws.Range("E121").Value = 638
$ws.Range("F121").Value = 45833.42523148148
# Row 123
$ws.Range("C123").Value = 416
$ws.Range("D123").Value = 45833.77618177892
$ws.Range("E123").Value = 416
$ws.Range("F123").Value = 45833.54157407407
# Row 124
$ws.Range("C124").Value = 253
$ws.Range("D124").Value = 45833.77618177133
$ws.Range("E124").Value = 253
$ws.Range("F124").Value = 45833.42833333334
# Row 125
$ws.Range("C125").Value = 513
$ws.Range("D125").Value = 45833.77620283911
$ws.Range("E125").Value = 513
$ws.Range("F125").Value = 45833.68387731481
# Row 126
$ws.Range("C126").Value = 31
$ws.Range("D126").Value = 45834.28681886589
$ws.Range("E126").Value = 31
$ws.Range("F126").Value = 45833.8759375
# Row 141
$ws.Range("C141").Value = 357
$ws.Range("D141").Value = 45833.77620284037
$ws.Range("E141").Value = 357
$ws.Range("F141").Value = 45833.68996527778
# Row 145
$ws.Range("C145").Value = 45
$ws.Range("D145").Value = 45833.77620284045
$ws.Range("E145").Value = 45
$ws.Range("F145").Value = 45833.68996527778
# Row 161
$ws.Range("C161").Value = 76
$ws.Range("D161").Value = 45833.77620284128
$ws.Range("E161").Value = 76
$ws.Range("F161").Value = 45833.73408564815
# Row 163
$ws.Range("C163").Value = 19
$ws.Range("D163").Value = 45833.77618177913
$ws.Range("E163").Value = 19
$ws.Range("F163").Value = 45833.54157407407
# Row 206
$ws.Range("C206").Value = 2
$ws.Range("D206").Value = 45833.77620284163
$ws.Range("E206").Value = 2
$ws.Range("F206").Value = 45833.73585648148
# Row 209
$ws.Range("C209").Value = 23
$ws.Range("D209").Value = 45833.77618177226
$ws.Range("E209").Value = 23
$ws.Range("F209").Value = 45833.44604166667
# Row 217
$ws.Range("C217").Value = 51
$ws.Range("D217").Value = 45833.77620283873
$ws.Range("E217").Value = 51
$ws.Range("F217").Value = 45833.66048611111
# Row 272
$ws.Range("C272").Value = 221
$ws.Range("D272").Value = 45833.77618177247
$ws.Range("E272").Value = 221
$ws.Range("F272").Value = 45833.44604166667
# Row 273
$ws.Range("C273").Value = 37
$ws.Range("D273").Value = 45833.77620283919
$ws.Range("E273").Value = 37
$ws.Range("F273").Value = 45833.68387731481
# Row 274
$ws.Range("C274").Value = 302
$ws.Range("D274").Value = 45833.77618177416
$ws.Range("E274").Value = 302
$ws.Range("F274").Value = 45833.48585648148
# Row 281
$ws.Range("C281").Value = 37
$ws.Range("D281").Value = 45833.77618177698
$ws.Range("E281").Value = 37
$ws.Range("F281").Value = 45833.52773148148
# Row 283
$ws.Range("C283").Value = 171
$ws.Range("D283").Value = 45833.77618177031
$ws.Range("E283").Value = 171
$ws.Range("F283").Value = 45833.42523148148
# Row 291
$ws.Range("C291").Value = 557
$ws.Range("D291").Value = 45833.77620284171
$ws.Range("E291").Value = 557
$ws.Range("F291").Value = 45833.73585648148
# Row 292
$ws.Range("C292").Value = 181
$ws.Range("D292").Value = 45833.77620284052
$ws.Range("E292").Value = 181
$ws.Range("F292").Value = 45833.68996527778
# Row 295
$ws.Range("C295").Value = 40
$ws.Range("D295").Value = 45833.77618177933
$ws.Range("E295").Value = 40
$ws.Range("F295").Value = 45833.54157407407
# Row 307
$ws.Range("C307").Value = 15
$ws.Range("D307").Value = 45833.77618177206
$ws.Range("E307").Value = 15
$ws.Range("F307").Value = 45833.44583333333
# Row 309
$ws.Range("C309").Value = 981
$ws.Range("D309").Value = 45833.77618177718
$ws.Range("E309").Value = 981
$ws.Range("F309").Value = 45833.53011574074
# Row 328
$ws.Range("C328").Value = 174
$ws.Range("D328").Value = 45833.77620283881
$ws.Range("E328").Value = 174
$ws.Range("F328").Value = 45833.66048611111
# Row 334
$ws.Range("C334").Value = 94
$ws.Range("D334").Value = 45833.77618177437
$ws.Range("E334").Value = 94
$ws.Range("F334").Value = 45833.48585648148
# Row 346
$ws.Range("C346").Value = 95
$ws.Range("D346").Value = 45833.77618177954
$ws.Range("E346").Value = 95
$ws.Range("F346").Value = 45833.54157407407
# Row 351
$ws.Range("C351").Value = 1238
$ws.Range("D351").Value = 45833.77618177974
$ws.Range("E351").Value = 1238
$ws.Range("F351").Value = 45833.54157407407
# Row 358
$ws.Range("C358").Value = 37
$ws.Range("D358").Value = 45833.77618177677
$ws.Range("E358").Value = 37
$ws.Range("F358").Value = 45833.52721064815
# Row 361
$ws.Range("C361").Value = 449
$ws.Range("D361").Value = 45833.77620284179
$ws.Range("E361").Value = 449
$ws.Range("F361").Value = 45833.73585648148
# Row 363
$ws.Range("C363").Value = 518
$ws.Range("D363").Value = 45833.77618177994
$ws.Range("E363").Value = 518
$ws.Range("F363").Value = 45833.54157407407
# Row 371
$ws.Range("C371").Value = 71
$ws.Range("D371").Value = 45833.77618177456
$ws.Range("E371").Value = 71
$ws.Range("F371").Value = 45833.48585648148
# Row 379
$ws.Range("C379").Value = 2
$ws.Range("D379").Value = 45833.77620283888
$ws.Range("E379").Value = 2
$ws.Range("F379").Value = 45833.66048611111
# Row 386
$ws.Range("C386").Value = 94
$ws.Range("D386").Value = 45833.77618178015
$ws.Range("E386").Value = 94
$ws.Range("F386").Value = 45833.54157407407
# Row 390
$ws.Range("C390").Value = 7
$ws.Range("D390").Value = 45833.77618177357
$ws.Range("E390").Value = 7
$ws.Range("F390").Value = 45833.45938657408
# Row 394
$ws.Range("C394").Value = 330
$ws.Range("D394").Value = 45833.77620284007
$ws.Range("E394").Value = 330
$ws.Range("F394").Value = 45833.68471064815
# Row 403
$ws.Range("C403").Value = 61
$ws.Range("D403").Value = 45833.77620284185
$ws.Range("E403").Value = 61
$ws.Range("F403").Value = 45833.73585648148
# Row 404
$ws.Range("C404").Value = -2
$ws.Range("D404").Value = 45833.77620283926
$ws.Range("E404").Value = -2
$ws.Range("F404").Value = 45833.68387731481
# Row 422
$ws.Range("C422").Value = 43
$ws.Range("D422").Value = 45833.77618177052
$ws.Range("E422").Value = 43
$ws.Range("F422").Value = 45833.42523148148
# Row 469
$ws.Range("C469").Value = 3046
$ws.Range("D469").Value = 45833.7762028406
$ws.Range("E469").Value = 3046
$ws.Range("F469").Value = 45833.68996527778
# Row 480
$ws.Range("C480").Value = 183
$ws.Range("D480").Value = 45833.77618177747
$ws.Range("E480").Value = 183
$ws.Range("F480").Value = 45833.53011574074
# Row 494
$ws.Range("C494").Value = 9
$ws.Range("D494").Value = 45833.77620284332
$ws.Range("E494").Value = 9
$ws.Range("F494").Value = 45833.769375
# Row 510
$ws.Range("C510").Value = 229
$ws.Range("D510").Value = 45833.77620283933
$ws.Range("E510").Value = 229
$ws.Range("F510").Value = 45833.68387731481
# Row 534
$ws.Range("C534").Value = 1228
$ws.Range("D534").Value = 45833.77620283941
$ws.Range("E534").Value = 1228
$ws.Range("F534").Value = 45833.68387731481
# Row 535
$ws.Range("C535").Value = 93
$ws.Range("D535").Value = 45833.77620283949
$ws.Range("E535").Value = 93
$ws.Range("F535").Value = 45833.68387731481
# Row 545
$ws.Range("C545").Value = 9
$ws.Range("D545").Value = 45833.77618177326
$ws.Range("E545").Value = 9
$ws.Range("F545").Value = 45833.45819444444
# Row 569
$ws.Range("C569").Value = 3
$ws.Range("D569").Value = 45833.77620283781
$ws.Range("E569").Value = 3
$ws.Range("F569").Value = 45833.63289351852
# Row 581
$ws.Range("C581").Value = 1
$ws.Range("D581").Value = 45833.7762028408
$ws.Range("E581").Value = 1
$ws.Range("F581").Value = 45833.68996527778
# Row 603
$ws.Range("C603").Value = 3
$ws.Range("D603").Value = 45833.77618177637
$ws.Range("E603").Value = 3
$ws.Range("F603").Value = 45833.515625
# Row 631
$ws.Range("C631").Value = 47
$ws.Range("D631").Value = 45833.77618177476
$ws.Range("E631").Value = 47
$ws.Range("F631").Value = 45833.48585648148
# Row 657
$ws.Range("C657").Value = 1870
$ws.Range("D657").Value = 45833.77620284192
$ws.Range("E657").Value = 1870
$ws.Range("F657").Value = 45833.73585648148
# Row 660
$ws.Range("C660").Value = 385
$ws.Range("D660").Value = 45833.77620284199
$ws.Range("E660").Value = 385
$ws.Range("F660").Value = 45833.73585648148
# Row 683
$ws.Range("C683").Value = 30
$ws.Range("D683").Value = 45833.77618177496
$ws.Range("E683").Value = 30
$ws.Range("F683").Value = 45833.48585648148
# Row 691
$ws.Range("C691").Value = 132
$ws.Range("D691").Value = 45833.77618177767
$ws.Range("E691").Value = 132
$ws.Range("F691").Value = 45833.53011574074
# Row 720
$ws.Range("C720").Value = 520
$ws.Range("D720").Value = 45833.77618177788
$ws.Range("E720").Value = 520
$ws.Range("F720").Value = 45833.53011574074
# Row 726
$ws.Range("C726").Value = 7
$ws.Range("D726").Value = 45833.77618178036
$ws.Range("E726").Value = 7
$ws.Range("F726").Value = 45833.54157407407
# Row 732
$ws.Range("C732").Value = 232
$ws.Range("D732").Value = 45833.77618177808
$ws.Range("E732").Value = 232
$ws.Range("F732").Value = 45833.53011574074
# Row 739
$ws.Range("C739").Value = 18
$ws.Range("D739").Value = 45833.77618178056
$ws.Range("E739").Value = 18
$ws.Range("F739").Value = 45833.54157407407
# Row 746
$ws.Range("C746").Value = 53
$ws.Range("D746").Value = 45833.77618177153
$ws.Range("E746").Value = 53
$ws.Range("F746").Value = 45833.42833333334
# Row 763
$ws.Range("C763").Value = -29
$ws.Range("D763").Value = 45833.77618176786
$ws.Range("E763").Value = -29
$ws.Range("F763").Value = 45833.38887731481
# Row 772
$ws.Range("C772").Value = 276
$ws.Range("D772").Value = 45833.77620284207
$ws.Range("E772").Value = 276
$ws.Range("F772").Value = 45833.73585648148
# Row 797
$ws.Range("C797").Value = 10
$ws.Range("D797").Value = 45833.77620284354
$ws.Range("E797").Value = 10
$ws.Range("F797").Value = 45833.79240740741
# Row 798
$ws.Range("C798").Value = 52
$ws.Range("D798").Value = 45833.77620284362
$ws.Range("E798").Value = 52
$ws.Range("F798").Value = 45833.79266203703
# Row 810
$ws.Range("C810").Value = 0
$ws.Range("D810").Value = 45833.77620283803
$ws.Range("E810").Value = 0
$ws.Range("F810").Value = 45833.65263888889
# Row 826
$ws.Range("C826").Value = 1
$ws.Range("D826").Value = 45833.77620284015
$ws.Range("E826").Value = 1
$ws.Range("F826").Value = 45833.68471064815
# Row 844
$ws.Range("C844").Value = -11
$ws.Range("D844").Value = 45833.77620283956
$ws.Range("E844").Value = -11
$ws.Range("F844").Value = 45833.68387731481
# Row 853
$ws.Range("C853").Value = 34
$ws.Range("D853").Value = 45833.77620284214
$ws.Range("E853").Value = 34
$ws.Range("F853").Value = 45833.73585648148
# Row 854
$ws.Range("C854").Value = 241
$ws.Range("D854").Value = 45834.28681886503
$ws.Range("E854").Value = 241
$ws.Range("F854").Value = 45833.87168981481
# Row 883
$ws.Range("C883").Value = 376
$ws.Range("D883").Value = 45833.77620284105
$ws.Range("E883").Value = 376
$ws.Range("F883").Value = 45833.69550925926
# Row 886
$ws.Range("D886").Value = 45833.7762028431
$ws.Range("F886").Value = 45833.76693287037
# Row 888
$ws.Range("C888").Value = 452
$ws.Range("D888").Value = 45833.77620284224
$ws.Range("E888").Value = 452
$ws.Range("F888").Value = 45833.73585648148
# Row 889
$ws.Range("C889").Value = 21
$ws.Range("D889").Value = 45833.77618177516
$ws.Range("E889").Value = 21
$ws.Range("F889").Value = 45833.48585648148
# Row 963
$ws.Range("C963").Value = 1618
$ws.Range("D963").Value = 45833.77618176898
$ws.Range("E963").Value = 1618
$ws.Range("F963").Value = 45833.39194444445
# Row 970
$ws.Range("C970").Value = 55
$ws.Range("D970").Value = 45833.77620284022
$ws.Range("E970").Value = 55
$ws.Range("F970").Value = 45833.68471064815
# Row 1002
$ws.Range("C1002").Value = 122
$ws.Range("D1002").Value = 45833.77618178075
$ws.Range("E1002").Value = 122
$ws.Range("F1002").Value = 45833.54157407407
# Row 1018
$ws.Range("C1018").Value = 6
$ws.Range("D1018").Value = 45834.28681886329
$ws.Range("E1018").Value = 6
$ws.Range("F1018").Value = 45833.85578703704
# Row 1024
$ws.Range("C1024").Value = 167
$ws.Range("D1024").Value = 45833.77620284279
$ws.Range("E1024").Value = 167
$ws.Range("F1024").Value = 45833.73675925926
# Row 1025
$ws.Range("C1025").Value = 240
$ws.Range("D1025").Value = 45833.77620283711
$ws.Range("E1025").Value = 240
$ws.Range("F1025").Value = 45833.54157407407
# Row 1054
$ws.Range("C1054").Value = 8
$ws.Range("D1054").Value = 45833.77620283797
$ws.Range("E1054").Value = 8
$ws.Range("F1054").Value = 45833.65202546296
# Row 1062
$ws.Range("C1062").Value = 258
$ws.Range("D1062").Value = 45833.77620284113
$ws.Range("E1062").Value = 258
$ws.Range("F1062").Value = 45833.69792824074
# Row 1094
$ws.Range("C1094").Value = 4
$ws.Range("D1094").Value = 45833.77618177071
$ws.Range("E1094").Value = 4
$ws.Range("F1094").Value = 45833.42523148148
# Row 1098
$ws.Range("C1098").Value = 39
$ws.Range("D1098").Value = 45833.7762028403
$ws.Range("E1098").Value = 39
$ws.Range("F1098").Value = 45833.68471064815
# Row 1110
$ws.Range("C1110").Value = 115
$ws.Range("D1110").Value = 45833.77618177092
$ws.Range("E1110").Value = 115
$ws.Range("F1110").Value = 45833.42523148148
# Row 1122
$ws.Range("C1122").Value = 18
$ws.Range("D1122").Value = 45833.77618176814
$ws.Range("E1122").Value = 18
$ws.Range("F1122").Value = 45833.38887731481
# Row 1126
$ws.Range("C1126").Value = 677
$ws.Range("D1126").Value = 45833.77620284229
$ws.Range("E1126").Value = 677
$ws.Range("F1126").Value = 45833.73585648148
# Row 1133
$ws.Range("C1133").Value = 72
$ws.Range("D1133").Value = 45833.77618177536
$ws.Range("E1133").Value = 72
$ws.Range("F1133").Value = 45833.48585648148
# Row 1147
$ws.Range("C1147").Value = 42
$ws.Range("D1147").Value = 45833.77620284236
$ws.Range("E1147").Value = 42
$ws.Range("F1147").Value = 45833.73585648148
# Row 1149
$ws.Range("C1149").Value = 6
$ws.Range("D1149").Value = 45833.77618177555
$ws.Range("E1149").Value = 6
$ws.Range("F1149").Value = 45833.48585648148
# Row 1150
$ws.Range("C1150").Value = 35
$ws.Range("D1150").Value = 45834.28681886363
$ws.Range("E1150").Value = 35
$ws.Range("F1150").Value = 45833.85730324074
# Row 1156
$ws.Range("C1156").Value = 15
$ws.Range("D1156").Value = 45833.77618177657
$ws.Range("E1156").Value = 15
$ws.Range("F1156").Value = 45833.515625
# Row 1174
$ws.Range("C1174").Value = 3
$ws.Range("D1174").Value = 45833.77618177576
$ws.Range("E1174").Value = 3
$ws.Range("F1174").Value = 45833.48585648148
# Row 1186
$ws.Range("C1186").Value = -1
$ws.Range("D1186").Value = 45833.77620284091
$ws.Range("E1186").Value = -1
$ws.Range("F1186").Value = 45833.68996527778
# Row 1228
$ws.Range("F1228").Value = 45833.77162037037
# Row 1250
$ws.Range("C1250").Value = 100
$ws.Range("D1250").Value = 45834.28681886417
$ws.Range("E1250").Value = 100
$ws.Range("F1250").Value = 45833.86076388889
# Row 1253
$ws.Range("C1253").Value = 279
$ws.Range("D1253").Value = 45833.77620284244
$ws.Range("E1253").Value = 279
$ws.Range("F1253").Value = 45833.73585648148
# Row 1255
$ws.Range("C1255").Value = 0
$ws.Range("D1255").Value = 45833.77618176919
$ws.Range("E1255").Value = 0
$ws.Range("F1255").Value = 45833.40459490741
# Row 1261
$ws.Range("C1261").Value = 253
$ws.Range("D1261").Value = 45834.28681886605
$ws.Range("E1261").Value = 253
$ws.Range("F1261").Value = 45833.91274305555
# Row 1324
$ws.Range("C1324").Value = 7
$ws.Range("D1324").Value = 45833.77620283811
$ws.Range("E1324").Value = 7
$ws.Range("F1324").Value = 45833.65372685185
# Row 1342
$ws.Range("C1342").Value = 542
$ws.Range("D1342").Value = 45833.77620283842
$ws.Range("E1342").Value = 542
$ws.Range("F1342").Value = 45833.66032407407
# Row 1385
$ws.Range("C1385").Value = 326
$ws.Range("D1385").Value = 45834.28681886399
$ws.Range("E1385").Value = 326
$ws.Range("F1385").Value = 45833.85784722222
# Row 1413
$ws.Range("C1413").Value = 64
$ws.Range("D1413").Value = 45833.7762028385
$ws.Range("E1413").Value = 64
$ws.Range("F1413").Value = 45833.66032407407
# Row 1422
$ws.Range("F1422").Value = 45833.76759259259
# Row 1439
$ws.Range("C1439").Value = 32
$ws.Range("D1439").Value = 45833.77618177827
$ws.Range("E1439").Value = 32
$ws.Range("F1439").Value = 45833.53011574074
# Row 1446
$ws.Range("C1446").Value = 46
$ws.Range("D1446").Value = 45833.77618177596
$ws.Range("E1446").Value = 46
$ws.Range("F1446").Value = 45833.48585648148
# Row 1448
$ws.Range("C1448").Value = 118
$ws.Range("D1448").Value = 45833.77620284251
$ws.Range("E1448").Value = 118
$ws.Range("F1448").Value = 45833.73585648148
# Row 1506
$ws.Range("C1506").Value = 10
$ws.Range("D1506").Value = 45833.77620284294
$ws.Range("E1506").Value = 10
$ws.Range("F1506").Value = 45833.76515046296
# Row 1507
$ws.Range("C1507").Value = 242
$ws.Range("D1507").Value = 45834.28681886305
$ws.Range("E1507").Value = 242
$ws.Range("F1507").Value = 45833.85539351852
# Row 1576
$ws.Range("C1576").Value = 43
$ws.Range("D1576").Value = 45833.77620283971
$ws.Range("E1576").Value = 43
$ws.Range("F1576").Value = 45833.68387731481
# Row 1594
$ws.Range("C1594").Value = 52
$ws.Range("D1594").Value = 45833.77620283723
$ws.Range("E1594").Value = 52
$ws.Range("F1594").Value = 45833.54157407407
# Row 1597
$ws.Range("C1597").Value = 4785
$ws.Range("D1597").Value = 45833.77620284258
$ws.Range("E1597").Value = 4785
$ws.Range("F1597").Value = 45833.73585648148
# Row 1600
$ws.Range("C1600").Value = 34
$ws.Range("D1600").Value = 45833.77620283866
$ws.Range("E1600").Value = 34
$ws.Range("F1600").Value = 45833.66032407407
# Row 1635
$ws.Range("C1635").Value = 2588
$ws.Range("D1635").Value = 45833.77620284265
$ws.Range("E1635").Value = 2588
$ws.Range("F1635").Value = 45833.73585648148
# Row 1651
$ws.Range("C1651").Value = 16
$ws.Range("D1651").Value = 45833.77620284302
$ws.Range("E1651").Value = 16
$ws.Range("F1651").Value = 45833.76623842592
# Row 1656
$ws.Range("C1656").Value = 15
$ws.Range("D1656").Value = 45834.28681886452
$ws.Range("E1656").Value = 15
$ws.Range("F1656").Value = 45833.87039351852
# Row 1657
$ws.Range("C1657").Value = 48
$ws.Range("D1657").Value = 45834.28681886521
$ws.Range("E1657").Value = 48
$ws.Range("F1657").Value = 45833.87199074074
# Row 1685
$ws.Range("F1685").Value = 45833.76810185185
# Row 1697
$ws.Range("C1697").Value = 2
$ws.Range("D1697").Value = 45833.77618177616
$ws.Range("E1697").Value = 2
$ws.Range("F1697").Value = 45833.48585648148
# Row 1758
$ws.Range("C1758").Value = 55
$ws.Range("D1758").Value = 45834.28681886347
$ws.Range("E1758").Value = 55
$ws.Range("F1758").Value = 45833.85653935185
# Row 1810
$ws.Range("C1810").Value = 15
$ws.Range("D1810").Value = 45833.77618177173
$ws.Range("E1810").Value = 15
$ws.Range("F1810").Value = 45833.43788194445
# Row 1824
$ws.Range("C1824").Value = 11
$ws.Range("D1824").Value = 45833.77618176959
$ws.Range("E1824").Value = 11
$ws.Range("F1824").Value = 45833.42410879629
# Row 1856
$ws.Range("C1856").Value = 3
$ws.Range("D1856").Value = 45833.77620283733
$ws.Range("E1856").Value = 3
$ws.Range("F1856").Value = 45833.54157407407
# Row 1861
$ws.Range("C1861").Value = 1
$ws.Range("D1861").Value = 45833.77620283978
$ws.Range("E1861").Value = 1
$ws.Range("F1861").Value = 45833.68387731481
# Row 1862
$ws.Range("F1862").Value = 45833.77251157408
# Row 1878
$ws.Range("C1878").Value = 43
$ws.Range("D1878").Value = 45833.77618177266
$ws.Range("E1878").Value = 43
$ws.Range("F1878").Value = 45833.44604166667
# Row 1912
$ws.Range("C1912").Value = 239
$ws.Range("D1912").Value = 45834.28681886538
$ws.Range("E1912").Value = 239
$ws.Range("F1912").Value = 45833.87310185185
# Row 2023
$ws.Range("C2023").Value = 102
$ws.Range("D2023").Value = 45834.2868188647
$ws.Range("E2023").Value = 102
$ws.Range("F2023").Value = 45833.87069444444
# Row 2024
$ws.Range("C2024").Value = 74
$ws.Range("D2024").Value = 45834.28681886487
$ws.Range("E2024").Value = 74
$ws.Range("F2024").Value = 45833.8713425926
# Row 2033
$ws.Range("C2033").Value = 65
$ws.Range("D2033").Value = 45834.28681886572
$ws.Range("E2033").Value = 65
$ws.Range("F2033").Value = 45833.87530092592
# Row 2062
$ws.Range("C2062").Value = -9
$ws.Range("D2062").Value = 45833.77620283788
$ws.Range("E2062").Value = -9
$ws.Range("F2062").Value = 45833.63289351852
# Row 2104
$ws.Range("C2104").Value = 18
$ws.Range("D2104").Value = 45833.77618176939
$ws.Range("E2104").Value = 18
$ws.Range("F2104").Value = 45833.40459490741
# Row 2144
$ws.Range("C2144").Value = 75
$ws.Range("D2144").Value = 45833.77618177286
$ws.Range("E2144").Value = 75
$ws.Range("F2144").Value = 45833.44604166667
# Row 2221
$ws.Range("C2221").Value = 18
$ws.Range("D2221").Value = 45833.77620283985
$ws.Range("E2221").Value = 18
$ws.Range("F2221").Value = 45833.68387731481
# Row 2241
$ws.Range("C2241").Value = 66
$ws.Range("D2241").Value = 45833.77620284271
$ws.Range("E2241").Value = 66
$ws.Range("F2241").Value = 45833.73585648148
# Row 2327
$ws.Range("C2327").Value = 14
$ws.Range("D2327").Value = 45833.7762028412
$ws.Range("E2327").Value = 14
$ws.Range("F2327").Value = 45833.69792824074
# Row 2379
$ws.Range("D2379").Value = 45833.77620283819
$ws.Range("F2379").Value = 45833.6553587963
$ws.Range("F2379").NumberFormat = $ws.Range("D2379").NumberFormat
# Row 2416
$ws.Range("C2416").Value = 96
$ws.Range("D2416").Value = 45833.77618177112
$ws.Range("E2416").Value = 96
$ws.Range("F2416").Value = 45833.42523148148
# Row 2465
$ws.Range("C2465").Value = 11
$ws.Range("D2465").Value = 45833.77618177306
$ws.Range("E2465").Value = 11
$ws.Range("F2465").Value = 45833.44604166667
# Row 2481
$ws.Range("C2481").Value = 203
$ws.Range("D2481").Value = 45833.77620283741
$ws.Range("E2481").Value = 203
$ws.Range("F2481").Value = 45833.54157407407
# Row 2482
$ws.Range("C2482").Value = 834
$ws.Range("D2482").Value = 45833.77620283749
$ws.Range("E2482").Value = 834
$ws.Range("F2482").Value = 45833.54157407407
# Row 2483
$ws.Range("C2483").Value = 1137
$ws.Range("D2483").Value = 45833.77620283757
$ws.Range("E2483").Value = 1137
$ws.Range("F2483").Value = 45833.54157407407
# Row 2484
$ws.Range("C2484").Value = 494
$ws.Range("D2484").Value = 45833.77620283764
$ws.Range("E2484").Value = 494
$ws.Range("F2484").Value = 45833.54157407407
# Row 2582
$ws.Range("F2582").Value = 45833.65579861111
$ws.Range("F2582").NumberFormat = $ws.Range("D2582").NumberFormat

# --- Append new rows 2587-2589 ---
# Row 2587
$ws.Range("A2587").Value = 43873951
$ws.Range("B2587").Value = 1
$ws.Range("C2587").Value = -1
$ws.Range("D2587").Value = 45833.7762028428
$ws.Range("D2587").NumberFormat = $ws.Range("D2586").NumberFormat
$ws.Range("E2587").Value = -1
$ws.Range("F2587").Value = 45833.68996527778
$ws.Range("F2587").NumberFormat = $ws.Range("D2587").NumberFormat
$ws.Range("G2587").Value = 0
$ws.Range("H2587").Value = "Consistente"
# Row 2588
$ws.Range("A2588").Value = 43874050
$ws.Range("B2588").Value = 1
$ws.Range("C2588").Value = 0
$ws.Range("D2588").Value = 45833.77620283766
$ws.Range("D2588").NumberFormat = $ws.Range("D2587").NumberFormat
$ws.Range("E2588").Value = 0
$ws.Range("G2588").Value = 0
$ws.Range("H2588").Value = "Consistente"
# Row 2589
$ws.Range("A2589").Value = 43883152
$ws.Range("B2589").Value = 1
$ws.Range("C2589").Value = 144
$ws.Range("D2589").Value = 45834.28681886436
$ws.Range("D2589").NumberFormat = $ws.Range("D2588").NumberFormat
$ws.Range("E2589").Value = 144
$ws.Range("F2589").Value = 45833.86863425926
$ws.Range("F2589").NumberFormat = $ws.Range("D2589").NumberFormat
$ws.Range("G2589").Value = 0
$ws.Range("H2589").Value = "Consistente"
